$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")
# (equivalently $wb.ActiveSheet, as this workbook has a single sheet)

$ws.Range('A1').Value = 'Datos actualizados a 20 de Junio de 2020 a las 13:33'
$ws.Range('B4').Value = 2297642
$ws.Range('C4').Value = 452
$ws.Range('E4').Value = 1220158
$ws.Range('B7').Value = 396765
$ws.Range('C7').Value = 953
$ws.Range('E7').Value = 169447
$ws.Range('G7').Value = 2
$ws.Range('H7').Value = 12972
$ws.Range('B13').Value = 202584
$ws.Range('C13').Value = 2322
$ws.Range('D13').Value = 161384
$ws.Range('E13').Value = 31693
$ws.Range('G13').Value = 115
$ws.Range('H13').Value = 9507
$ws.Range('B27').Value = 57936
$ws.Range('C27').Value = 603
$ws.Range('D27').Value = 36749
$ws.Range('E27').Value = 20844
$ws.Range('G27').Value = 6
$ws.Range('H27').Value = 343
$ws.Range('B33').Value = 44533
$ws.Range('C33').Value = 388
$ws.Range('D33').Value = 31754
$ws.Range('E33').Value = 12478
$ws.Range('G33').Value = 1
$ws.Range('H33').Value = 301
$ws.Range('B40').Value = 31243
$ws.Range('C40').Value = 8
$ws.Range('E40').Value = 387
$ws.Range('A41').Value = 'Filipinas'
$ws.Range('B41').Value = 29400
$ws.Range('C41').Value = 941
$ws.Range('D41').Value = 7650
$ws.Range('E41').Value = 20600
$ws.Range('G41').Value = 20
$ws.Range('H41').Value = 1150
$ws.Range('A42').Value = 'Oman'
$ws.Range('B42').Value = 28566
$ws.Range('C42').Value = 896
$ws.Range('D42').Value = 14780
$ws.Range('E42').Value = 13658
$ws.Range('G42').Value = 3
$ws.Range('H42').Value = 128
$ws.Range('A70').Value = 'Nepal'
$ws.Range('B70').Value = 8605
$ws.Range('C70').Value = 331
$ws.Range('D70').Value = 1578
$ws.Range('E70').Value = 7005
$ws.Range('H70').Value = 22
$ws.Range('A71').Value = 'Malasia'
$ws.Range('B71').Value = 8556
$ws.Range('C71').Value = 21
$ws.Range('D71').Value = 8146
$ws.Range('E71').Value = 289
$ws.Range('H71').Value = 121
$ws.Range('A72').Value = 'Sudan'
$ws.Range('B72').Value = 8316
$ws.Range('D72').Value = 3086
$ws.Range('E72').Value = 4724
$ws.Range('H72').Value = 506
$ws.Range('A77').Value = 'Senegal'
$ws.Range('B77').Value = 5783
$ws.Range('C77').Value = 144
$ws.Range('D77').Value = 3859
$ws.Range('E77').Value = 1842
$ws.Range('H77').Value = 82
$ws.Range('A78').Value = 'Consejo Danes para los Refugiados'
$ws.Range('B78').Value = 5672
$ws.Range('C78').Value = 195
$ws.Range('D78').Value = 807
$ws.Range('E78').Value = 4740
$ws.Range('G78').Value = 3
$ws.Range('H78').Value = 125
$ws.Range('B119').Value = 1503
$ws.Range('C119').Value = 60
$ws.Range('D119').Value = 618
$ws.Range('E119').Value = 872
$ws.Range('B124').Value = 1129
$ws.Range('C124').Value = 1
$ws.Range('D124').Value = 1077
$ws.Range('E124').Value = 48
$ws.Range('A138').Value = 'Estado de Palestina'
$ws.Range('B138').Value = 759
$ws.Range('C138').Value = 84
$ws.Range('D138').Value = 437
$ws.Range('E138').Value = 319
$ws.Range('H138').Value = 3
$ws.Range('A139').Value = 'Crucero'
$ws.Range('B139').Value = 712
$ws.Range('C139').Value = 0
$ws.Range('D139').Value = 651
$ws.Range('E139').Value = 48
$ws.Range('H139').Value = 13
$ws.Range('B143').Value = 664
$ws.Range('C143').Value = 1
$ws.Range('D143').Value = 616
$ws.Range('E143').Value = 39
$ws.Range('A202').Value = 'Dominica'
$ws.Range('A203').Value = 'Fiyi'
$ws.Range('A208').Value = 'Islas Turcas y Caicos'
$ws.Range('D208').Value = 11
$ws.Range('H208').Value = 1
$ws.Range('A209').Value = 'Santa Sede'
$ws.Range('D209').Value = 12
$ws.Range('H209').Value = 0
